# Trade #25 closed at 2026-02-16 22:54:34 - base_strategy UP +0.000%
# Append a new trade row (row 26) to both the "All Trades" sheet and the
# "base_strategy" sheet, mirroring the existing rows (e.g. row 25).

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 26

    # Trade # (number)
    $ws.Cells.Item($row, 1).Value = 25

    # Date - force literal text so Excel does not auto-convert it to a date serial.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-16"

    # Time (stays text on its own, no special handling needed)
    $ws.Cells.Item($row, 3).Value = "22:54:34"

    # Strategy
    $ws.Cells.Item($row, 4).Value = "base_strategy"

    # Side
    $ws.Cells.Item($row, 5).Value = "UP"

    # Entry Price
    $ws.Cells.Item($row, 6).Value = 49.999998

    # Exit Price - left blank/empty (trade is still OPEN)
    $ws.Cells.Item($row, 7).Value = ""

    # Status
    $ws.Cells.Item($row, 8).Value = "OPEN"

    # P&L %
    $ws.Cells.Item($row, 9).Value = 0

    # P&L $
    $ws.Cells.Item($row, 10).Value = 0

    # Capital After
    $ws.Cells.Item($row, 11).Value = 100

    # Entry Slippage (bps)
    $ws.Cells.Item($row, 12).Value = 0

    # Exit Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0

    # Confidence
    $ws.Cells.Item($row, 14).Value = 0.6

    # Entry Reason
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"

    # Exit Reason - left blank/empty (trade is still OPEN)
    $ws.Cells.Item($row, 16).Value = ""

    # Duration (min)
    $ws.Cells.Item($row, 17).Value = 0
}
